$d = $word.ActiveDocument

# Phase 1: replace each old value with a unique temporary placeholder
# to avoid cascading replacements when a new value matches another old value.
$d.Content.Find.Execute("79÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0@@", 2) | Out-Null
$d.Content.Find.Execute("68÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH1@@", 2) | Out-Null
$d.Content.Find.Execute("97÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH2@@", 2) | Out-Null
$d.Content.Find.Execute("25÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH3@@", 2) | Out-Null
$d.Content.Find.Execute("48÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH4@@", 2) | Out-Null
$d.Content.Find.Execute("53÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH5@@", 2) | Out-Null
$d.Content.Find.Execute("69÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH6@@", 2) | Out-Null
$d.Content.Find.Execute("54÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH7@@", 2) | Out-Null
$d.Content.Find.Execute("65÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH8@@", 2) | Out-Null
$d.Content.Find.Execute("62÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH9@@", 2) | Out-Null
$d.Content.Find.Execute("82÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH10@@", 2) | Out-Null
$d.Content.Find.Execute("84÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH11@@", 2) | Out-Null
$d.Content.Find.Execute("32÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH12@@", 2) | Out-Null
$d.Content.Find.Execute("41÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH13@@", 2) | Out-Null
$d.Content.Find.Execute("48÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH14@@", 2) | Out-Null
$d.Content.Find.Execute("83÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH15@@", 2) | Out-Null
$d.Content.Find.Execute("49÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH16@@", 2) | Out-Null
$d.Content.Find.Execute("92÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH17@@", 2) | Out-Null
$d.Content.Find.Execute("62÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH18@@", 2) | Out-Null
$d.Content.Find.Execute("93÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH19@@", 2) | Out-Null
$d.Content.Find.Execute("89÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH20@@", 2) | Out-Null
$d.Content.Find.Execute("98÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH21@@", 2) | Out-Null
$d.Content.Find.Execute("51÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH22@@", 2) | Out-Null
$d.Content.Find.Execute("21÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH23@@", 2) | Out-Null
$d.Content.Find.Execute("27÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH24@@", 2) | Out-Null

# Phase 2: replace each placeholder with its final new value
$d.Content.Find.Execute("@@PH0@@", $true, $false, $false, $false, $false, $true, 1, $false, "50÷3=", 2) | Out-Null
$d.Content.Find.Execute("@@PH1@@", $true, $false, $false, $false, $false, $true, 1, $false, "48÷5=", 2) | Out-Null
$d.Content.Find.Execute("@@PH2@@", $true, $false, $false, $false, $false, $true, 1, $false, "87÷8=", 2) | Out-Null
$d.Content.Find.Execute("@@PH3@@", $true, $false, $false, $false, $false, $true, 1, $false, "31÷2=", 2) | Out-Null
$d.Content.Find.Execute("@@PH4@@", $true, $false, $false, $false, $false, $true, 1, $false, "69÷4=", 2) | Out-Null
$d.Content.Find.Execute("@@PH5@@", $true, $false, $false, $false, $false, $true, 1, $false, "62÷9=", 2) | Out-Null
$d.Content.Find.Execute("@@PH6@@", $true, $false, $false, $false, $false, $true, 1, $false, "45÷8=", 2) | Out-Null
$d.Content.Find.Execute("@@PH7@@", $true, $false, $false, $false, $false, $true, 1, $false, "33÷5=", 2) | Out-Null
$d.Content.Find.Execute("@@PH8@@", $true, $false, $false, $false, $false, $true, 1, $false, "25÷5=", 2) | Out-Null
$d.Content.Find.Execute("@@PH9@@", $true, $false, $false, $false, $false, $true, 1, $false, "22÷7=", 2) | Out-Null
$d.Content.Find.Execute("@@PH10@@", $true, $false, $false, $false, $false, $true, 1, $false, "37÷5=", 2) | Out-Null
$d.Content.Find.Execute("@@PH11@@", $true, $false, $false, $false, $false, $true, 1, $false, "51÷2=", 2) | Out-Null
$d.Content.Find.Execute("@@PH12@@", $true, $false, $false, $false, $false, $true, 1, $false, "49÷5=", 2) | Out-Null
$d.Content.Find.Execute("@@PH13@@", $true, $false, $false, $false, $false, $true, 1, $false, "84÷8=", 2) | Out-Null
$d.Content.Find.Execute("@@PH14@@", $true, $false, $false, $false, $false, $true, 1, $false, "84÷9=", 2) | Out-Null
$d.Content.Find.Execute("@@PH15@@", $true, $false, $false, $false, $false, $true, 1, $false, "38÷5=", 2) | Out-Null
$d.Content.Find.Execute("@@PH16@@", $true, $false, $false, $false, $false, $true, 1, $false, "66÷2=", 2) | Out-Null
$d.Content.Find.Execute("@@PH17@@", $true, $false, $false, $false, $false, $true, 1, $false, "30÷6=", 2) | Out-Null
$d.Content.Find.Execute("@@PH18@@", $true, $false, $false, $false, $false, $true, 1, $false, "11÷2=", 2) | Out-Null
$d.Content.Find.Execute("@@PH19@@", $true, $false, $false, $false, $false, $true, 1, $false, "80÷4=", 2) | Out-Null
$d.Content.Find.Execute("@@PH20@@", $true, $false, $false, $false, $false, $true, 1, $false, "20÷9=", 2) | Out-Null
$d.Content.Find.Execute("@@PH21@@", $true, $false, $false, $false, $false, $true, 1, $false, "21÷6=", 2) | Out-Null
$d.Content.Find.Execute("@@PH22@@", $true, $false, $false, $false, $false, $true, 1, $false, "34÷4=", 2) | Out-Null
$d.Content.Find.Execute("@@PH23@@", $true, $false, $false, $false, $false, $true, 1, $false, "29÷7=", 2) | Out-Null
$d.Content.Find.Execute("@@PH24@@", $true, $false, $false, $false, $false, $true, 1, $false, "79÷8=", 2) | Out-Null
